$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": a new daily column ("04-dec") is inserted right
#     before the old "01-oct." column (EG), shifting every column from EG
#     onward one place to the right (EG -> EH, ... FK -> FL). ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Range("EG1").EntireColumn.Insert()
$wsPrix.Range("EG1").Value = "04-dec"
$wsPrix.Range("EG2:EG25").Value = "-"

# --- Sheet "Gaz": one more day of data is appended as row 167. ---
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force text formatting first so the ISO-looking date string isn't
# auto-converted into a date serial number, then drop the number format
# again so the new cell ends up plain/unstyled like its neighbours.
$wsGaz.Range("A167").NumberFormat = "@"
$wsGaz.Range("A167").Value = "2025-12-02"
$wsGaz.Range("A167").ClearFormats()
$wsGaz.Range("B167").Value = 26.895

# --- Sheet "CO2": matching new row 167. ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A167").NumberFormat = "@"
$wsCO2.Range("A167").Value = "2025-12-02"
$wsCO2.Range("A167").ClearFormats()
$wsCO2.Range("B167").Value = 81.65000000000001
